$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add two new cage records (rows 36-37) below the existing table.
$ws.Range("A36").Value = "85B"
$ws.Range("B36").Value = 320
$ws.Range("C36").Value = 320
$ws.Range("D36").Value = 320
$ws.Range("E36").Value = "Metal"

$ws.Range("A37").Value = "942B"
$ws.Range("B37").Value = 250
$ws.Range("C37").Value = 250
$ws.Range("D37").Value = 250
$ws.Range("E37").Value = "Metal"
